# Calculator Keyboard Layout - add 8051 side keyboard code and bitmaps
# for summation (Sigma) and product (Pi); add bracket keys; add Enter key.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- H6 "=" -> "Enter" (smaller Consolas 12pt so the word fits) -----------
$ws.Range("H6").Value = "Enter"
$ws.Range("H6").Font.Name = "Consolas"
$ws.Range("H6").Font.Size = 12

# --- L4/M4 ( ) -> [ ] and W4/X4 ( ) -> { } ---------------------------------
# These four keys move from the plain (white) style to the highlighted
# (yellow) style already used elsewhere (e.g. R4). Copy that formatting
# over first, then set the new glyph values.
$fmtSrc = $ws.Range("R4")
$fmtSrc.Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("W4").PasteSpecial(-4122)
$ws.Range("X4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("L4").Value = "["
$ws.Range("M4").Value = "]"
$ws.Range("W4").Value = "{"
$ws.Range("X4").Value = "}"

# --- U5 restyle (yellow -> white, matches existing M5/X5 look) ------------
$fmtSrc = $ws.Range("M5")
$fmtSrc.Copy()
$ws.Range("U5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- AD6 restyle (white -> yellow, matches existing R4/S6 look) -----------
$fmtSrc = $ws.Range("R4")
$fmtSrc.Copy()
$ws.Range("AD6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Q6 "EE" -> "x" (multiply glyph, matches F6/AB6) ----------------------
$fmtSrc = $ws.Range("F6")
$fmtSrc.Copy()
$ws.Range("Q6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Q6").Value = "×"

# --- Update the stored selection to match the saved workbook state --------
[void]$ws.Range("W8").Select()
